# "Generate Report for Handback" - mark the zh-cn and de-de handback rows as
# synced, stamp the handback datetime, and record the generated target /
# handback file names (with a hyperlink on the target file, like the
# existing "source file" column already has).

$wb = $excel.ActiveWorkbook

$ghBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/80e7b6856553e831bd26e12f9a107371d49986d5/e2e/d120e4fd-5ce2-49e9-aba3-09a2b49ec8aa.md"
$mdName = "d120e4fd-5ce2-49e9-aba3-09a2b49ec8aa.md"

# ---- zh-cn sheet -------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = "Handed back: in sync with en-US"

$zh.Range("I2").Value = $mdName
$zh.Hyperlinks.Add($zh.Range("I2"), $ghBase, "", "", $mdName) | Out-Null

$zh.Range("J2").Value = "d120e4fd-5ce2-49e9-aba3-09a2b49ec8aa.db2273b2c4a71febf59e71415699cafeaac3f6ce.zh-cn.xlf"
$zh.Range("K2").Value = "2016-09-03 13:04:17"

$zh.Columns.Item(3).ColumnWidth = 29.1666666666667
$zh.Columns.Item(9).ColumnWidth = 39.1666666666667
$zh.Columns.Item(10).ColumnWidth = 39.1666666666667

# ---- de-de sheet --------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = "Handed back: in sync with en-US"

$de.Range("I2").Value = $mdName
$de.Hyperlinks.Add($de.Range("I2"), $ghBase, "", "", $mdName) | Out-Null

$de.Range("J2").Value = "d120e4fd-5ce2-49e9-aba3-09a2b49ec8aa.db2273b2c4a71febf59e71415699cafeaac3f6ce.de-de.xlf"
$de.Range("K2").Value = "2016-09-03 13:04:24"

$de.Columns.Item(3).ColumnWidth = 29.1666666666667
$de.Columns.Item(9).ColumnWidth = 39.1666666666667
$de.Columns.Item(10).ColumnWidth = 39.1666666666667

# ---- Overview sheet -------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("E2").Value = "Handed back: in sync with en-US"
$ov.Range("F2").Value = "Handed back: in sync with en-US"

$ov.Columns.Item(5).ColumnWidth = 29.1666666666667
$ov.Columns.Item(6).ColumnWidth = 29.1666666666667
